$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44432
$ws.Cells.Item(2, 10).Value = 30
$ws.Cells.Item(2, 11).Value = 7000
$ws.Cells.Item(2, 12).Value = 7000
$ws.Cells.Item(2, 13).Value = 7000
$ws.Cells.Item(2, 16).Value = 583

$ws.Cells.Item(3, 4).Value = 44195
$ws.Cells.Item(3, 10).Value = 55
$ws.Cells.Item(3, 12).Value = 5000
$ws.Cells.Item(3, 13).Value = 5000
$ws.Cells.Item(3, 16).Value = 417

$ws.Cells.Item(4, 4).Value = 44428
$ws.Cells.Item(4, 10).Value = 10
$ws.Cells.Item(4, 11).Value = 7000
$ws.Cells.Item(4, 12).Value = 7000
$ws.Cells.Item(4, 13).Value = 7000
$ws.Cells.Item(4, 16).Value = 583

$ws.Cells.Item(5, 4).Value = 44497
$ws.Cells.Item(5, 10).Value = 40

$ws.Cells.Item(6, 4).Value = 44495
$ws.Cells.Item(6, 10).Value = 30

$ws.Cells.Item(7, 4).Value = 44438
$ws.Cells.Item(7, 10).Value = 30

$ws.Cells.Item(8, 4).Value = 44452
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = 7000
$ws.Cells.Item(8, 16).Value = 583

$ws.Cells.Item(9, 4).Value = 44203
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 5000
$ws.Cells.Item(9, 16).Value = 417

$ws.Cells.Item(10, 4).Value = 44424
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 7000
$ws.Cells.Item(10, 16).Value = 583

$ws.Cells.Item(11, 4).Value = 44498
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 7000
$ws.Cells.Item(11, 12).Value = 7000
$ws.Cells.Item(11, 13).Value = 7000
$ws.Cells.Item(11, 16).Value = 583

$ws.Cells.Item(13, 4).Value = 44410
$ws.Cells.Item(13, 10).Value = 40
$ws.Cells.Item(13, 11).Value = 7000
$ws.Cells.Item(13, 12).Value = 7000
$ws.Cells.Item(13, 13).Value = 7000
$ws.Cells.Item(13, 16).Value = 583

$ws.Cells.Item(14, 4).Value = 44441
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 6000
$ws.Cells.Item(14, 16).Value = 500

$ws.Cells.Item(15, 4).Value = 44413
$ws.Cells.Item(15, 10).Value = 40
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 7000
$ws.Cells.Item(15, 13).Value = 7000
$ws.Cells.Item(15, 16).Value = 583

$ws.Cells.Item(16, 4).Value = 44448

$ws.Cells.Item(17, 4).Value = 44165
$ws.Cells.Item(17, 10).Value = 130
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 13).Value = 5615
$ws.Cells.Item(17, 16).Value = 468

$ws.Cells.Item(18, 4).Value = 44467
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 6000
$ws.Cells.Item(18, 13).Value = 6000
$ws.Cells.Item(18, 16).Value = 500

$ws.Cells.Item(19, 4).Value = 44455
$ws.Cells.Item(19, 10).Value = 20

$ws.Cells.Item(20, 4).Value = 44483
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 7000
$ws.Cells.Item(20, 12).Value = 8000
$ws.Cells.Item(20, 13).Value = 7600
$ws.Cells.Item(20, 16).Value = 633

$ws.Cells.Item(21, 4).Value = 44435
$ws.Cells.Item(21, 10).Value = 30

$ws.Cells.Item(22, 4).Value = 44427
$ws.Cells.Item(22, 10).Value = 20
$ws.Cells.Item(22, 12).Value = 7000
$ws.Cells.Item(22, 13).Value = 7000
$ws.Cells.Item(22, 16).Value = 583

$ws.Cells.Item(23, 4).Value = 44259

$ws.Cells.Item(24, 4).Value = 44327
$ws.Cells.Item(24, 10).Value = 30
$ws.Cells.Item(24, 11).Value = 6000
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 6000
$ws.Cells.Item(24, 16).Value = 500

$ws.Cells.Item(25, 4).Value = 44162
$ws.Cells.Item(25, 10).Value = 50

$ws.Cells.Item(26, 4).Value = 44326
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 6000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 6000
$ws.Cells.Item(26, 16).Value = 500

$ws.Cells.Item(27, 4).Value = 44474

$ws.Cells.Item(28, 4).Value = 44369
$ws.Cells.Item(28, 10).Value = 20
$ws.Cells.Item(28, 11).Value = 4000
$ws.Cells.Item(28, 12).Value = 4000
$ws.Cells.Item(28, 13).Value = 4000
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 333

$ws.Cells.Item(29, 4).Value = 44196
$ws.Cells.Item(29, 11).Value = 5000
$ws.Cells.Item(29, 12).Value = 5000
$ws.Cells.Item(29, 13).Value = 5000
$ws.Cells.Item(29, 16).Value = 417

$ws.Cells.Item(30, 4).Value = 44454

$ws.Cells.Item(32, 4).Value = 44490
$ws.Cells.Item(32, 10).Value = 65
$ws.Cells.Item(32, 11).Value = 6000
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = 6000
$ws.Cells.Item(32, 16).Value = 500

$ws.Cells.Item(33, 4).Value = 44329
$ws.Cells.Item(33, 11).Value = 5000
$ws.Cells.Item(33, 12).Value = 6000
$ws.Cells.Item(33, 13).Value = 5500
$ws.Cells.Item(33, 16).Value = 458

$ws.Cells.Item(34, 4).Value = 44453
$ws.Cells.Item(34, 10).Value = 20
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 12).Value = 7000
$ws.Cells.Item(34, 13).Value = 7000
$ws.Cells.Item(34, 16).Value = 583

$ws.Cells.Item(35, 4).Value = 44466
$ws.Cells.Item(35, 10).Value = 50
$ws.Cells.Item(35, 11).Value = 6000
$ws.Cells.Item(35, 13).Value = 6400
$ws.Cells.Item(35, 16).Value = 533

$ws.Cells.Item(36, 4).Value = 44442
$ws.Cells.Item(36, 10).Value = 20
$ws.Cells.Item(36, 11).Value = 6000
$ws.Cells.Item(36, 12).Value = 6000
$ws.Cells.Item(36, 13).Value = 6000
$ws.Cells.Item(36, 16).Value = 500

$ws.Cells.Item(37, 4).Value = 44476
$ws.Cells.Item(37, 10).Value = 30
$ws.Cells.Item(37, 11).Value = 8000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 8000
$ws.Cells.Item(37, 16).Value = 667

$ws.Cells.Item(38, 4).Value = 44372
$ws.Cells.Item(38, 10).Value = 40
$ws.Cells.Item(38, 11).Value = 6000
$ws.Cells.Item(38, 12).Value = 6000
$ws.Cells.Item(38, 13).Value = 6000
$ws.Cells.Item(38, 16).Value = 500

$ws.Cells.Item(39, 4).Value = 44166
$ws.Cells.Item(39, 10).Value = 55
$ws.Cells.Item(39, 11).Value = 6000
$ws.Cells.Item(39, 12).Value = 6000
$ws.Cells.Item(39, 13).Value = 6000
$ws.Cells.Item(39, 16).Value = 500

$ws.Cells.Item(40, 4).Value = 44477
$ws.Cells.Item(40, 11).Value = 8000
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = 8000
$ws.Cells.Item(40, 16).Value = 667

$ws.Cells.Item(41, 4).Value = 44211
$ws.Cells.Item(41, 10).Value = 65
$ws.Cells.Item(41, 11).Value = 5000
$ws.Cells.Item(41, 12).Value = 5000
$ws.Cells.Item(41, 13).Value = 5000
$ws.Cells.Item(41, 16).Value = 417

$ws.Cells.Item(42, 4).Value = 44186
$ws.Cells.Item(42, 10).Value = 50
$ws.Cells.Item(42, 11).Value = 5000
$ws.Cells.Item(42, 12).Value = 5000
$ws.Cells.Item(42, 13).Value = 5000
$ws.Cells.Item(42, 16).Value = 417

$ws.Cells.Item(43, 4).Value = 44487
$ws.Cells.Item(43, 10).Value = 55

$ws.Cells.Item(44, 4).Value = 44425
$ws.Cells.Item(44, 11).Value = 7000
$ws.Cells.Item(44, 12).Value = 7000
$ws.Cells.Item(44, 13).Value = 7000
$ws.Cells.Item(44, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(44, 16).Value = 583

$ws.Cells.Item(45, 4).Value = 44449
$ws.Cells.Item(45, 10).Value = 65
$ws.Cells.Item(45, 11).Value = 7000
$ws.Cells.Item(45, 12).Value = 7000
$ws.Cells.Item(45, 13).Value = 7000
$ws.Cells.Item(45, 16).Value = 583

$ws.Cells.Item(46, 4).Value = 44484
$ws.Cells.Item(46, 10).Value = 30
$ws.Cells.Item(46, 12).Value = 8000
$ws.Cells.Item(46, 13).Value = 7333
$ws.Cells.Item(46, 16).Value = 611

$ws.Cells.Item(47, 4).Value = 44494
$ws.Cells.Item(47, 10).Value = 30
$ws.Cells.Item(47, 12).Value = 6000
$ws.Cells.Item(47, 13).Value = 6000
$ws.Cells.Item(47, 16).Value = 500

$ws.Cells.Item(48, 4).Value = 44301
$ws.Cells.Item(48, 10).Value = 50

$ws.Cells.Item(49, 4).Value = 44302
$ws.Cells.Item(49, 10).Value = 20

$ws.Cells.Item(50, 4).Value = 44179
$ws.Cells.Item(50, 10).Value = 40
$ws.Cells.Item(50, 11).Value = 6000
$ws.Cells.Item(50, 12).Value = 6000
$ws.Cells.Item(50, 13).Value = 6000
$ws.Cells.Item(50, 16).Value = 500
